$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129; this shifts the existing rows 129-153
# down to 130-154 (and the sheet dimension grows from T153 to T154).
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with a new weekly price record.
# Most columns repeat the same market/product metadata as the surrounding
# rows; only the date, volume, prices and $/kg differ.
$ws.Range("A129").Value = 4
$ws.Range("B129").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C129").Value = "Los Lagos"
$ws.Range("D129").Value = 44504
$ws.Range("E129").Value = 10
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100108
$ws.Range("H129").Value = "Tropicales y subtropicales"
$ws.Range("I129").Value = 100108005
$ws.Range("J129").Value = "Piña"
$ws.Range("K129").Value = "Caramelo"
$ws.Range("L129").Value = "Segunda"
$ws.Range("M129").Value = 120
$ws.Range("N129").Value = 20000
$ws.Range("O129").Value = 21000
$ws.Range("P129").Value = 20500
$ws.Range("Q129").Value = "$/caja 14 unidades"
$ws.Range("R129").Value = "Ecuador"
$ws.Range("S129").Value = 1464
$ws.Range("T129").Value = 14

# Make sure the D129 cell keeps the date number format used by the rest
# of column D (style index 2 in the original workbook).
$ws.Range("D129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
